{"js": "// Find the table-cell text \"1 Change\" (the row whose neighboring cell is\n// \"PC_1\") and append \" (New description)\" right after it, producing\n// \"1 Change (New description)\". Using a search hit (instead of inserting\n// at the paragraph's \"End\") anchors the insertion immediately after the\n// existing \"...nge\" run so the new text inherits that run's formatting\n// (w:lang=\"en-US\"), matching the rest of the cell's text.\nconst body = context.document.body;\nconst results = body.search(\"1 Change\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '1 Change' text\");\n}\n\nresults.items[0].insertText(\" (New description)\", \"End\");\nawait context.sync();\n", "ps1": "# Locate the table-cell text \"1 Change\" (the row whose neighboring cell is\n# \"PC_1\") and append \" (New description)\" right after it, producing\n# \"1 Change (New description)\".\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = \"1 Change\"\n$r.Find.MatchCase = $true\n$r.Find.MatchWholeWord = $false\n$found = $r.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the '1 Change' text\"\n}\n\n# Collapse the found range to its end point (right after \"...Change\") and\n# insert the new text there so it inherits the adjacent run's formatting\n# (w:lang=\"en-US\"), matching the rest of the cell's text.\n$r.Collapse($wdCollapseEnd)\n$r.InsertAfter(\" (New description)\")\n"}
